$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new data row (row 16) that was added on July 22
$ws.Range("A16").Value = 44763
$ws.Range("B16").Value = 50
$ws.Range("C16").Value = 201
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 44
$ws.Range("F16").Value = 190
$ws.Range("G16").Value = 8

# Update the selected cell/view to reflect where the user left off
$ws.Range("L29").Select()
